# Scheduled-runner style refresh of cached market-price figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across the
# per-job-class Leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 990.0833
$ws.Range("I19").Value = 784.5714
$ws.Range("K19").Value = 784.5714
$ws.Range("M19").Value = -609.5714

$ws.Range("H43").Value = 5333.1665
$ws.Range("I43").Value = 4666.3335
$ws.Range("K43").Value = 4666.3335
$ws.Range("M43").Value = -4597.3335

$ws.Range("H58").Value = 2807.5
$ws.Range("I58").Value = 282.8
$ws.Range("J58").Value = 3955.0908
$ws.Range("K58").Value = 848.4000000000001
$ws.Range("L58").Value = 11865.2724
$ws.Range("M58").Value = -698.4000000000001
$ws.Range("N58").Value = -12165.2724

$ws.Range("H74").Value = 4518.625
$ws.Range("I74").Value = 4518.625
$ws.Range("K74").Value = 4518.625
$ws.Range("M74").Value = -3582.625

$ws.Range("H77").Value = 4518.625
$ws.Range("I77").Value = 4518.625
$ws.Range("K77").Value = 22593.125
$ws.Range("M77").Value = -17913.125

$ws.Range("H88").Value = 4415.6665
$ws.Range("I88").Value = 4770.25
$ws.Range("J88").Value = 4132
$ws.Range("K88").Value = 4770.25
$ws.Range("L88").Value = 4132
$ws.Range("M88").Value = -4364.25
$ws.Range("N88").Value = -4944

$ws.Range("H91").Value = 4415.6665
$ws.Range("I91").Value = 4770.25
$ws.Range("J91").Value = 4132
$ws.Range("K91").Value = 4770.25
$ws.Range("L91").Value = 4132
$ws.Range("M91").Value = -3366.25
$ws.Range("N91").Value = -6940

$ws.Range("H111").Value = 549.5
$ws.Range("I111").Value = 549.5
$ws.Range("K111").Value = 1648.5
$ws.Range("M111").Value = 1418.5

$ws.Range("H116").Value = 42499.668
$ws.Range("J116").Value = 59999.5
$ws.Range("L116").Value = 59999.5
$ws.Range("N116").Value = -66883.5

$ws.Range("H132").Value = 4424.5
$ws.Range("I132").Value = 4424.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13273.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10743.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2147.2334
$ws.Range("I32").Value = 1200.6296
$ws.Range("K32").Value = 1200.6296
$ws.Range("M32").Value = -913.6296

$ws.Range("H102").Value = 12903.308
$ws.Range("I102").Value = 10774.5
$ws.Range("J102").Value = 19999.334
$ws.Range("K102").Value = 10774.5
$ws.Range("L102").Value = 19999.334
$ws.Range("M102").Value = -9152.5
$ws.Range("N102").Value = -23243.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 11032.526
$ws.Range("I94").Value = 7663.077
$ws.Range("K94").Value = 7663.077
$ws.Range("M94").Value = -7212.077

$ws.Range("H134").Value = 4969.6665
$ws.Range("J134").Value = 2913.3333
$ws.Range("L134").Value = 8739.999899999999
$ws.Range("N134").Value = -13809.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 34336.973
$ws.Range("I22").Value = 37777.156
$ws.Range("J22").Value = 6815.5
$ws.Range("K22").Value = 37777.156
$ws.Range("L22").Value = 6815.5
$ws.Range("M22").Value = -37427.156
$ws.Range("N22").Value = -7515.5

$ws.Range("H31").Value = 8184.6875
$ws.Range("I31").Value = 3746.875
$ws.Range("J31").Value = 12622.5
$ws.Range("K31").Value = 3746.875
$ws.Range("L31").Value = 12622.5
$ws.Range("M31").Value = -3451.875
$ws.Range("N31").Value = -13212.5

$ws.Range("H34").Value = 8184.6875
$ws.Range("I34").Value = 3746.875
$ws.Range("J34").Value = 12622.5
$ws.Range("K34").Value = 3746.875
$ws.Range("L34").Value = 12622.5
$ws.Range("M34").Value = -3544.875
$ws.Range("N34").Value = -13026.5

$ws.Range("H96").Value = 9147.267
$ws.Range("I96").Value = 300
$ws.Range("J96").Value = 9779.214
$ws.Range("K96").Value = 300
$ws.Range("L96").Value = 9779.214
$ws.Range("N96").Value = -15271.214
$ws.Range("M96").Value = 2446

$ws.Range("H114").Value = 28228
$ws.Range("J114").Value = 28228
$ws.Range("L114").Value = 28228
$ws.Range("N114").Value = -36906

$ws.Range("H132").Value = 4572.5713
$ws.Range("I132").Value = 3665.6667
$ws.Range("K132").Value = 10997.0001
$ws.Range("M132").Value = -8467.000100000001

$ws.Range("H134").Value = 3175.3333
$ws.Range("I134").Value = 2069.8235
$ws.Range("K134").Value = 6209.470499999999
$ws.Range("M134").Value = -3674.470499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1091.6818
$ws.Range("I11").Value = 197.92857
$ws.Range("J11").Value = 1508.7667
$ws.Range("K11").Value = 593.78571
$ws.Range("L11").Value = 4526.300099999999
$ws.Range("M11").Value = -453.78571
$ws.Range("N11").Value = -4806.300099999999

$ws.Range("H122").Value = 1900
$ws.Range("J122").Value = 1900
$ws.Range("L122").Value = 17100
$ws.Range("N122").Value = -22000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6644.636
$ws.Range("I43").Value = 3309.6
$ws.Range("K43").Value = 3309.6
$ws.Range("M43").Value = -3158.6

$ws.Range("H46").Value = 4112.2856
$ws.Range("I46").Value = 4112.2856
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 4112.2856
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -3956.2856
$ws.Range("N46").ClearContents()

$ws.Range("H57").Value = 18258.334
$ws.Range("I57").Value = 9888
$ws.Range("J57").Value = 34999
$ws.Range("K57").Value = 9888
$ws.Range("L57").Value = 34999
$ws.Range("M57").Value = -9068
$ws.Range("N57").Value = -36639

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2148.75
$ws.Range("I22").Value = 1531.6666
$ws.Range("K22").Value = 1531.6666
$ws.Range("M22").Value = -1236.6666

$ws.Range("H27").Value = 2148.75
$ws.Range("I27").Value = 1531.6666
$ws.Range("K27").Value = 1531.6666
$ws.Range("M27").Value = -1424.6666

$ws.Range("H55").Value = 4150.125
$ws.Range("I55").Value = 2706.6875
$ws.Range("J55").Value = 7037
$ws.Range("K55").Value = 2706.6875
$ws.Range("L55").Value = 7037
$ws.Range("M55").Value = -2533.6875
$ws.Range("N55").Value = -7383

$ws.Range("H100").Value = 11499.375
$ws.Range("I100").Value = 9999.5
$ws.Range("K100").Value = 9999.5
$ws.Range("M100").Value = -9458.5

$ws.Range("H132").Value = 39998
$ws.Range("I132").Value = 39998
$ws.Range("K132").Value = 119994
$ws.Range("M132").Value = -117464

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2066.5557
$ws.Range("J81").Value = 2700
$ws.Range("L81").Value = 5400
$ws.Range("N81").Value = -7522

$ws.Range("H84").Value = 2066.5557
$ws.Range("J84").Value = 2700
$ws.Range("L84").Value = 27000
$ws.Range("N84").Value = -37608

$ws.Range("H93").Value = 55000
$ws.Range("J93").Value = 55000
$ws.Range("L93").Value = 55000
$ws.Range("N93").Value = -59992

$ws.Range("H136").Value = 4835
$ws.Range("I136").Value = 3823.6667
$ws.Range("K136").Value = 11471.0001
$ws.Range("M136").Value = -8921.000100000001
